# Dataframe ST.xlsx — "Add files via upload" commit
#
# What changed (per the OOXML diff):
#   1. Sheet3's small lookup/source table (A20:B36) got refreshed numbers
#      for the 23-Nov pull (column B).
#   2. Sheet1 gained one more trailing date column (CO) for "23-nov",
#      mirroring the existing CN ("22-nov") column: a text header in CO1
#      and, for every data row, a plain numeric value in CO equal to the
#      freshly recomputed VLOOKUP result (same value CB/CC now resolve to).
#   3. Everything else (CB/CC formulas, Sheet3 C2:C18 formulas) updates by
#      itself once the source table changes and the workbook recalculates.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- 1. Refresh the Sheet3 source table (A20:B36) with the new pull ------
$newValues = @{
    "3D QUESO 92GX27"                      = 8.4232659615610537
    "CHEETOS QUESO 85GX24X1"               = 3.9927386838936654
    "DORITOS QUESO 129GX19"                = 7.1287121915406404
    "DORITOS QUESO 70X40G"                 = 9.7470478670968514
    "DORITOS QUESO 77GX26"                 = 6.8906842034833167
    "LAYS CEBOLLA CARAMELIZADA 85GX25"     = 0
    "LAYS CLASICAS 145GRX18"               = 15.641916725891482
    "LAYS CLASICAS 249GRX14"               = 9.4634220040227603
    "LAYS CLASICAS 40GX68"                 = 6.3568830725813941
    "LAYS CLASICAS 94GRX25"                = 0.30438357978260722
    "LAYS ONDAS FH 30GX72"                 = 10.773000000033225
    "LAYS ONDAS FH 70GX28"                 = 9.9957550000833084
    "LAYS QSO Y CEBOLLA 34GX72"            = 23.045471303355239
    "PEHUAMAR ACANALADA 520GX9"            = 8.7384188426844194
    "PEHUAMAR MAICITOS 285GX10"            = 10.919534699567107
    "PEHUAMAR PAPA LISA 520GX9"            = 9.2372786663899991
    "QUAKER AVENA INSTANT FORTIF 18X280G"  = 37.719044642736968
}

for ($r = 20; $r -le 36; $r++) {
    $name = $ws3.Cells.Item($r, 1).Value()
    if ($newValues.ContainsKey($name)) {
        $ws3.Cells.Item($r, 2).Value = $newValues[$name]
    }
}

# Let CB/CC (Sheet1) and C2:C18 (Sheet3) formulas pick up the new source.
$excel.Calculate()

# --- 2. Add the new "23-nov" column (CO) on Sheet1 -----------------------
$ws1.Cells.Item(1, 93).Value = "23-nov"

for ($r = 2; $r -le 18; $r++) {
    # CO mirrors CB/CC (the recalculated VLOOKUP result) as a plain value,
    # same as how CN already holds a plain copy of that lookup.
    $cb = $ws1.Cells.Item($r, 80).Value()
    $ws1.Cells.Item($r, 93).Value = $cb
    $ws1.Cells.Item($r, 93).NumberFormat = "0"
}

$excel.Calculate()

# --- 3. Leave the selection where the author left it (CQ7) ---------------
$ws1.Activate() | Out-Null
$ws1.Range("CQ7").Select() | Out-Null
